$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.406.17"
$ws.Range("E2").Value = "  -3.95%  "
$ws.Range("D3").Value = "3.005.35"
$ws.Range("E3").Value = "  -2.74%  "
$ws.Range("E4").Value = "  -0.11%  "
$ws.Range("D5").Value = "'548.45"
$ws.Range("E5").Value = "  +0.92%  "
$ws.Range("D6").Value = "'134.50"
$ws.Range("E6").Value = "  -3.88%  "
$ws.Range("E7").Value = "  -0.10%  "
$ws.Range("D8").Value = "3.001.31"
$ws.Range("E8").Value = "  -2.68%  "
$ws.Range("E9").Value = "  -0.37%  "
$ws.Range("E10").Value = "  -5.02%  "
$ws.Range("E11").Value = "  -7.97%  "
$ws.Range("D12").Value = "'0.451"
$ws.Range("E12").Value = "  -1.55%  "
$ws.Range("D13").Value = "'34.57"
$ws.Range("E13").Value = "  -0.54%  "
$ws.Range("E14").Value = "  -2.30%  "
$ws.Range("D15").Value = "3.495.82"
$ws.Range("E15").Value = "  -2.87%  "
$ws.Range("D16").Value = "61.520.97"
$ws.Range("E16").Value = "  -3.93%  "
$ws.Range("E17").Value = "  -2.25%  "
$ws.Range("D18").Value = "3.002.73"
$ws.Range("E18").Value = "  -3.16%  "
$ws.Range("D19").Value = "'6.67"
$ws.Range("E19").Value = "  +0.20%  "
$ws.Range("D20").Value = "'473.27"
$ws.Range("E20").Value = "  -1.73%  "
$ws.Range("D21").Value = "'13.29"
$ws.Range("E21").Value = "  -0.92%  "
$ws.Range("D22").Value = "'0.676"
$ws.Range("E22").Value = "  -3.46%  "
$ws.Range("D23").Value = "'7.05"
$ws.Range("E23").Value = "  -0.88%  "
$ws.Range("D24").Value = "'80.08"
$ws.Range("E24").Value = "  +1.15%  "
$ws.Range("D25").Value = "'12.12"
$ws.Range("E25").Value = "  -2.22%  "
$ws.Range("D26").Value = "'1.00"
$ws.Range("E26").Value = "  +0.03%  "
$ws.Range("D27").Value = "'2.72"
$ws.Range("E27").Value = "  -0.57%  "
$ws.Range("D28").Value = "'7.82"
$ws.Range("E28").Value = "  -3.13%  "
$ws.Range("D29").Value = "'0.999"
$ws.Range("E29").Value = "  -0.16%  "
$ws.Range("E30").Value = "  +0.37%  "
$ws.Range("D31").Value = "'25.70"
$ws.Range("E31").Value = "  -2.25%  "
$ws.Range("E32").Value = "  -1.85%  "
$ws.Range("D33").Value = "'5.55"
$ws.Range("E33").Value = "  +2.98%  "
$ws.Range("E34").Value = "  -2.63%  "
$ws.Range("D35").Value = "'55.38"
$ws.Range("E35").Value = "  -3.34%  "
$ws.Range("D36").Value = "'5.91"
$ws.Range("E36").Value = "  -1.73%  "
$ws.Range("D37").Value = "'454.29"
$ws.Range("E37").Value = "  -8.30%  "
$ws.Range("D38").Value = "3.190.35"
$ws.Range("E38").Value = "  -2.73%  "
$ws.Range("D39").Value = "'0.0799"
$ws.Range("E39").Value = "  +0.00%  "
$ws.Range("D40").Value = "'0.0383"
$ws.Range("E40").Value = "  -5.47%  "
$ws.Range("D41").Value = "'0.118"
$ws.Range("E41").Value = "  -2.19%  "
$ws.Range("D42").Value = "'8.16"
$ws.Range("E42").Value = "  +0.72%  "
$ws.Range("D43").Value = "'2.44"
$ws.Range("E43").Value = "  -10.07%  "
$ws.Range("E44").Value = "  -0.01%  "
$ws.Range("D45").Value = "'26.21"
$ws.Range("E45").Value = "  +3.87%  "
$ws.Range("D46").Value = "'0.245"
$ws.Range("E46").Value = "  -3.54%  "
$ws.Range("D47").Value = "'1.99"
$ws.Range("E47").Value = "  -2.95%  "
$ws.Range("E48").Value = "  -0.60%  "
$ws.Range("D49").Value = "'118.09"
$ws.Range("E49").Value = "  -4.66%  "
$ws.Range("E50").Value = "  +6.96%  "
$ws.Range("D51").Value = "0.0₃0492"
$ws.Range("E51").Value = "  -7.61%  "
